# Update simulated run-time / error metrics for the "other language" low-input run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each data row (2-11), a list of [columnLetter, newValue] pairs covering
# run_time (C), max_er (E) and iter 0..19 (F:Y). Column D (num_deaths) is untouched.
$updates = @(
    @(2, @(@("C", 1.166001081466675), @("E", 1905.949421026427), @("F", 0.1336276190416302), @("G", 0.0982439249395764), @("H", 0.07543290281764277), @("I", 0.06864783755953557), @("J", 0.05848606940282208), @("K", 0.05273187124861078), @("L", 0.05023793305181231), @("M", 0.04740259355738934), @("N", 0.0440156755675446), @("O", 0.04251324011342084), @("P", 0.04120135785159825), @("Q", 0.03998900694595772), @("R", 0.03992414247785133), @("S", 0.03897165338332452), @("T", 0.03869947486989641), @("U", 0.03818965030896423), @("V", 0.03789584297132012), @("W", 0.0375685889534671), @("X", 0.03732017189733961), @("Y", 0.03715301015646057))),
    @(3, @(@("C", 1.359998464584351), @("E", 1824.659199608774), @("F", 0.1325212854253535), @("G", 0.09830096683322269), @("H", 0.06905253754975017), @("I", 0.06064223679925009), @("J", 0.05154669795054532), @("K", 0.04969908852991612), @("L", 0.04607070191828877), @("M", 0.04306320666139202), @("N", 0.04040280696085036), @("O", 0.03972213360615853), @("P", 0.03865794468531609), @("Q", 0.03757865875740376), @("R", 0.03749459972822405), @("S", 0.03708521544648542), @("T", 0.03636889663631319), @("U", 0.03621611334283644), @("V", 0.03596073159007885), @("W", 0.03562376097726265), @("X", 0.03562376097726265), @("Y", 0.03556840545046341))),
    @(4, @(@("C", 1.186997890472412), @("E", 1878.163667730181), @("F", 0.126696430494691), @("G", 0.09022642224688786), @("H", 0.0774839113390384), @("I", 0.06303023381947856), @("J", 0.05734017632120101), @("K", 0.05175459986853218), @("L", 0.0469567592677923), @("M", 0.04470477622716963), @("N", 0.04288315751981753), @("O", 0.0413883879447148), @("P", 0.03983158175519631), @("Q", 0.03875564487236852), @("R", 0.03866005991358947), @("S", 0.03808626685977744), @("T", 0.03764471726127983), @("U", 0.03742136020812947), @("V", 0.03685361075405776), @("W", 0.03685361075405776), @("X", 0.03663226764188224), @("Y", 0.03661137753860001))),
    @(5, @(@("C", 1.242001533508301), @("E", 1761.080355877262), @("F", 0.136867930071421), @("G", 0.09643128711759497), @("H", 0.0737848153196406), @("I", 0.06137087636325438), @("J", 0.05280470713389059), @("K", 0.04732123733656717), @("L", 0.04388290857979677), @("M", 0.04105685091981452), @("N", 0.03963930701941125), @("O", 0.03914497456160382), @("P", 0.03799613069637295), @("Q", 0.03655282990645944), @("R", 0.0357974268741478), @("S", 0.0351933291806543), @("T", 0.03484771283892917), @("U", 0.03484771283892917), @("V", 0.03458144848604628), @("W", 0.03457016193906156), @("X", 0.03443648975844797), @("Y", 0.03432905177148659))),
    @(6, @(@("C", 1.170998811721802), @("E", 1826.666673185373), @("F", 0.1368864481900249), @("G", 0.09905298246506937), @("H", 0.07943362386793017), @("I", 0.06504037414700746), @("J", 0.05826452093533624), @("K", 0.05269596057952335), @("L", 0.04847365907333451), @("M", 0.04521765387350313), @("N", 0.04292661011786137), @("O", 0.0413257609343315), @("P", 0.04046453873161566), @("Q", 0.03874300544439693), @("R", 0.03789335346122086), @("S", 0.03692675991207781), @("T", 0.03643257144655158), @("U", 0.03619799231435805), @("V", 0.03578720789809184), @("W", 0.03574094596271262), @("X", 0.0356789070684336), @("Y", 0.03560753748899362))),
    @(7, @(@("C", 1.052002668380737), @("E", 1844.537716492403), @("F", 0.1355283607395215), @("G", 0.09690146705006931), @("H", 0.08108314503030156), @("I", 0.06965851348985574), @("J", 0.05871170763278336), @("K", 0.05305234051353432), @("L", 0.04911119081491815), @("M", 0.04500451950167123), @("N", 0.04448308821927766), @("O", 0.04210919094789722), @("P", 0.04040512136016771), @("Q", 0.03899300976856405), @("R", 0.03855522361015514), @("S", 0.03732823562869549), @("T", 0.03689100279275073), @("U", 0.03674500932785378), @("V", 0.03648511593990093), @("W", 0.03627372517902573), @("X", 0.03601960374606182), @("Y", 0.03595590090628465))),
    @(8, @(@("C", 1.341996669769287), @("E", 1826.69292853496), @("F", 0.133713709466977), @("G", 0.1030438197670616), @("H", 0.07966442326432682), @("I", 0.06850593779019164), @("J", 0.05741519585426614), @("K", 0.05328651570010589), @("L", 0.04911999366269498), @("M", 0.04585073936838294), @("N", 0.04370039545059348), @("O", 0.04065636743945738), @("P", 0.03913423761535903), @("Q", 0.03816917561457289), @("R", 0.03764284110849648), @("S", 0.03673586775370491), @("T", 0.0364034550555528), @("U", 0.03611740091957354), @("V", 0.03584013067508935), @("W", 0.03578710764669066), @("X", 0.03569904997598621), @("Y", 0.03560804928918049))),
    @(9, @(@("C", 1.066000461578369), @("E", 1851.517649274663), @("F", 0.1374613389252656), @("G", 0.09371976525044966), @("H", 0.08224030088370353), @("I", 0.06850881005748262), @("J", 0.05855270088398399), @("K", 0.05203415667921735), @("L", 0.04838300138406228), @("M", 0.04567977807337065), @("N", 0.042606484294713), @("O", 0.04150007426121753), @("P", 0.03954062741801743), @("Q", 0.03841466265130356), @("R", 0.03770415153161494), @("S", 0.03733207457834017), @("T", 0.03678983764577225), @("U", 0.03659211755826047), @("V", 0.03652836407389177), @("W", 0.03626950721392352), @("X", 0.03622581547867144), @("Y", 0.03609196197416496))),
    @(10, @(@("C", 1.198999881744385), @("E", 1862.409047204985), @("F", 0.1293499576885841), @("G", 0.1033937403065161), @("H", 0.07932350838552298), @("I", 0.06501663577840935), @("J", 0.05872217903012848), @("K", 0.05194877695962238), @("L", 0.04672563930604243), @("M", 0.04492751804631682), @("N", 0.0433931164931188), @("O", 0.04075839205722648), @("P", 0.04046567955148522), @("Q", 0.03889803235334647), @("R", 0.03832190609137021), @("S", 0.03762098723454178), @("T", 0.03726417683206392), @("U", 0.03706618380021307), @("V", 0.03685038108494265), @("W", 0.03661669618059692), @("X", 0.03642635003873098), @("Y", 0.03630426992602308))),
    @(11, @(@("C", 1.094008445739746), @("E", 1834.292166806421), @("F", 0.1262916432714495), @("G", 0.08967207712067368), @("H", 0.07612048049073943), @("I", 0.06417723975115938), @("J", 0.05694102592608619), @("K", 0.05136793995655132), @("L", 0.04802779891832511), @("M", 0.0460594767022075), @("N", 0.04287652300610758), @("O", 0.04150719874883466), @("P", 0.03976349694591073), @("Q", 0.03853515690943396), @("R", 0.03685722792341475), @("S", 0.03685722792341475), @("T", 0.03684794093417776), @("U", 0.03647480103202949), @("V", 0.03614547556745465), @("W", 0.03604103317070813), @("X", 0.03579360454179947), @("Y", 0.03575618258881911)))
)

foreach ($rowUpdate in $updates) {
    $rowNum = $rowUpdate[0]
    $cellPairs = $rowUpdate[1]
    foreach ($pair in $cellPairs) {
        $colLetter = $pair[0]
        $newValue = $pair[1]
        $ws.Range("$colLetter$rowNum").Value = $newValue
    }
}